# Insert a new data row at row 137, shifting existing rows 137-276 down to 138-277.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(137).Insert()

# Populate the newly inserted row 137 with its data.
$row = 137
$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44539
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108002
$ws.Cells.Item($row, 10).Value = "Mango"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 684
$ws.Cells.Item($row, 14).Value = 6000
$ws.Cells.Item($row, 15).Value = 7000
$ws.Cells.Item($row, 16).Value = 6500
$ws.Cells.Item($row, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item($row, 18).Value = "Perú"
$ws.Cells.Item($row, 19).Value = 1625
$ws.Cells.Item($row, 20).Value = 4

# Ensure the new date cell keeps the same date-time number format as the other
# rows in column D (style index 2 in the original workbook).
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
